# "la inn unik emne funksjonalitet" - add "Unike emner" (unique topics) column
# functionality to the Statistikk sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistikk")

# --- 1. New column F: "Unike emner" header + values ------------------------
$ws.Range("F1").Value = "Unike emner"

$f2f31 = @(0,0,0,0,0,0,0,0,5,0,4,0,5,15,0,0,0,0,0,0,0,0,0,0,2,9,0,0,9,0)
for ($i = 0; $i -lt $f2f31.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 6).Value = $f2f31[$i]
}

# Give the new column a sensible custom width (closest the engine can store
# to the authored 12.5546875 "characters" width).
$ws.Columns.Item(6).ColumnWidth = 11.71

# --- 2. Nudge the two charts so they keep their visual position now that -----
#        a new column sits underneath/behind them.
$co1 = $ws.ChartObjects().Item(1)
$co1.Left = 411.1669921875
$co1.Top = 79.2
$co1.Width = 366.4658203125
$co1.Height = 216.0

$co2 = $ws.ChartObjects().Item(2)
$co2.Left = 779.8828125
$co2.Top = 79.2
$co2.Width = 356.5908203125
$co2.Height = 216.0

# --- 3. Make "Statistikk" the active tab/sheet and select F12 --------------
$ws.Activate()
$ws.Range("F12").Select()
